# Add 3-year RP (lead 3d) rows for 2025-10-30 and 2025-10-31 to both
# station hydrograph sheets (G5368 - Nia Pumping Station, G4945 -
# Talacogon Municipal Hall), extending each sheet from 30 to 32 data rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet "G5368" (Nia Pumping Station, primary) ----
$ws1 = $wb.Worksheets.Item("G5368")

$ws1.Cells.Item(31, 1).Value = "Philippines"
$ws1.Cells.Item(31, 2).Value = "philippines"
$ws1.Cells.Item(31, 3).Value = "Agusan River Basin"
$ws1.Cells.Item(31, 4).Value = "agusan"
$ws1.Cells.Item(31, 5).Value = "Nia Pumping Station"
$ws1.Cells.Item(31, 6).Value = "G5368"
$ws1.Cells.Item(31, 7).Value = "primary"
$ws1.Cells.Item(31, 8).Value = "'2025-10-30"
$ws1.Cells.Item(31, 8).Style = "Normal"
$ws1.Cells.Item(31, 9).Value = 3
$ws1.Cells.Item(31, 10).Value = 8.874999999999865
$ws1.Cells.Item(31, 11).Value = 125.5749999999995
$ws1.Cells.Item(31, 12).Value = 5
$ws1.Cells.Item(31, 13).Value = 4709.973879596918
$ws1.Cells.Item(31, 14).Value = "LOW"
$ws1.Cells.Item(31, 15).Value = 3193.342710267902
$ws1.Cells.Item(31, 16).Value = 4709.973879596918
$ws1.Cells.Item(31, 17).Value = 50
$ws1.Cells.Item(31, 18).Value = 0
$ws1.Cells.Item(31, 19).Value = 0
$ws1.Cells.Item(31, 20).Value = 982.4140625
$ws1.Cells.Item(31, 21).Value = 997.2823486328125
$ws1.Cells.Item(31, 22).Value = 735.6953125
$ws1.Cells.Item(31, 23).Value = 1445.9765625
$ws1.Cells.Item(31, 24).Value = 881.123046875
$ws1.Cells.Item(31, 25).Value = 1074.212890625
$ws1.Cells.Item(31, 26).Value = $false
$ws1.Cells.Item(31, 27).Value = -79.14183629009689

$ws1.Cells.Item(32, 1).Value = "Philippines"
$ws1.Cells.Item(32, 2).Value = "philippines"
$ws1.Cells.Item(32, 3).Value = "Agusan River Basin"
$ws1.Cells.Item(32, 4).Value = "agusan"
$ws1.Cells.Item(32, 5).Value = "Nia Pumping Station"
$ws1.Cells.Item(32, 6).Value = "G5368"
$ws1.Cells.Item(32, 7).Value = "primary"
$ws1.Cells.Item(32, 8).Value = "'2025-10-31"
$ws1.Cells.Item(32, 8).Style = "Normal"
$ws1.Cells.Item(32, 9).Value = 3
$ws1.Cells.Item(32, 10).Value = 8.874999999999865
$ws1.Cells.Item(32, 11).Value = 125.5749999999995
$ws1.Cells.Item(32, 12).Value = 5
$ws1.Cells.Item(32, 13).Value = 4709.973879596918
$ws1.Cells.Item(32, 14).Value = "LOW"
$ws1.Cells.Item(32, 15).Value = 3193.342710267902
$ws1.Cells.Item(32, 16).Value = 4709.973879596918
$ws1.Cells.Item(32, 17).Value = 50
$ws1.Cells.Item(32, 18).Value = 0
$ws1.Cells.Item(32, 19).Value = 0
$ws1.Cells.Item(32, 20).Value = 869.40234375
$ws1.Cells.Item(32, 21).Value = 878.986572265625
$ws1.Cells.Item(32, 22).Value = 734.5859375
$ws1.Cells.Item(32, 23).Value = 1094.109375
$ws1.Cells.Item(32, 24).Value = 800.96484375
$ws1.Cells.Item(32, 25).Value = 964.27734375
$ws1.Cells.Item(32, 26).Value = $false
$ws1.Cells.Item(32, 27).Value = -81.54124914543254

# ---- Sheet "G4945" (Talacogon Municipal Hall, secondary) ----
$ws2 = $wb.Worksheets.Item("G4945")

$ws2.Cells.Item(31, 1).Value = "Philippines"
$ws2.Cells.Item(31, 2).Value = "philippines"
$ws2.Cells.Item(31, 3).Value = "Agusan River Basin"
$ws2.Cells.Item(31, 4).Value = "agusan"
$ws2.Cells.Item(31, 5).Value = "Talacogon Municipal Hall"
$ws2.Cells.Item(31, 6).Value = "G4945"
$ws2.Cells.Item(31, 7).Value = "secondary"
$ws2.Cells.Item(31, 8).Value = "'2025-10-30"
$ws2.Cells.Item(31, 8).Style = "Normal"
$ws2.Cells.Item(31, 9).Value = 3
$ws2.Cells.Item(31, 10).Value = 8.424999999999859
$ws2.Cells.Item(31, 11).Value = 125.7749999999995
$ws2.Cells.Item(31, 12).Value = 5
$ws2.Cells.Item(31, 13).Value = 3363.250778297076
$ws2.Cells.Item(31, 14).Value = "LOW"
$ws2.Cells.Item(31, 15).Value = 2342.691130371584
$ws2.Cells.Item(31, 16).Value = 3363.250778297076
$ws2.Cells.Item(31, 17).Value = 50
$ws2.Cells.Item(31, 18).Value = 0
$ws2.Cells.Item(31, 19).Value = 0
$ws2.Cells.Item(31, 20).Value = 676.140625
$ws2.Cells.Item(31, 21).Value = 680.6256103515625
$ws2.Cells.Item(31, 22).Value = 488.8515625
$ws2.Cells.Item(31, 23).Value = 985.921875
$ws2.Cells.Item(31, 24).Value = 588.794921875
$ws2.Cells.Item(31, 25).Value = 745.98828125
$ws2.Cells.Item(31, 26).Value = $false
$ws2.Cells.Item(31, 27).Value = -79.89621739291317

$ws2.Cells.Item(32, 1).Value = "Philippines"
$ws2.Cells.Item(32, 2).Value = "philippines"
$ws2.Cells.Item(32, 3).Value = "Agusan River Basin"
$ws2.Cells.Item(32, 4).Value = "agusan"
$ws2.Cells.Item(32, 5).Value = "Talacogon Municipal Hall"
$ws2.Cells.Item(32, 6).Value = "G4945"
$ws2.Cells.Item(32, 7).Value = "secondary"
$ws2.Cells.Item(32, 8).Value = "'2025-10-31"
$ws2.Cells.Item(32, 8).Style = "Normal"
$ws2.Cells.Item(32, 9).Value = 3
$ws2.Cells.Item(32, 10).Value = 8.424999999999859
$ws2.Cells.Item(32, 11).Value = 125.7749999999995
$ws2.Cells.Item(32, 12).Value = 5
$ws2.Cells.Item(32, 13).Value = 3363.250778297076
$ws2.Cells.Item(32, 14).Value = "LOW"
$ws2.Cells.Item(32, 15).Value = 2342.691130371584
$ws2.Cells.Item(32, 16).Value = 3363.250778297076
$ws2.Cells.Item(32, 17).Value = 50
$ws2.Cells.Item(32, 18).Value = 0
$ws2.Cells.Item(32, 19).Value = 0
$ws2.Cells.Item(32, 20).Value = 555.4453125
$ws2.Cells.Item(32, 21).Value = 570.4659423828125
$ws2.Cells.Item(32, 22).Value = 459.1796875
$ws2.Cells.Item(32, 23).Value = 736.265625
$ws2.Cells.Item(32, 24).Value = 509.15625
$ws2.Cells.Item(32, 25).Value = 623.27734375
$ws2.Cells.Item(32, 26).Value = $false
$ws2.Cells.Item(32, 27).Value = -83.48486779266457
